$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (and its tab label is the sheet name itself)
$ws.Name = "Through 2021-11-05"

# Row 13 - "November (through 11-04/05)" row
$ws.Range("A13").Value = "November (through 11-05)"
$ws.Range("C13").Value = 5
$ws.Range("I13").Value = 19
$ws.Range("J13").Value = 0.05
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 0.1765
$ws.Range("O13").Value = 7
$ws.Range("R13").Value = 36
$ws.Range("U13").Value = 31

# Row 14 - "Total" row
$ws.Range("C14").Value = 231
$ws.Range("D14").Value = 0.1217
$ws.Range("I14").Value = 668
$ws.Range("J14").Value = 0.0849
$ws.Range("K14").Value = 69
$ws.Range("M14").Value = 0.1092
$ws.Range("O14").Value = 441
$ws.Range("P14").Value = 0.0982
$ws.Range("R14").Value = 1039
$ws.Range("S14").Value = 0.0494
$ws.Range("U14").Value = 1392
$ws.Range("V14").Value = 0.0563
